$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing empty (but styled) rows 6 and 7
$ws.Rows("6:7").Delete() | Out-Null

# Add the new "Observación" header in O5
$ws.Range("O5").Value = "Observación"

# Copy the formatting of the neighboring header cell (N5) onto O5,
# so the new header looks like the rest of the row 5 headers
$ws.Range("N5").Copy() | Out-Null
$ws.Range("O5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Give the new column a sensible width, matching the other header columns
$ws.Columns("O").ColumnWidth = 11.17

# Leave the selection on the newly added cell
$ws.Range("O5").Select() | Out-Null
